# Update cryptos list values (prices and 1h volume %) per the Apr 12 2023 refresh.
# Rows 38/39 (TheSandbox / Algorand) swap order with updated figures.
# Numeric-looking price strings are forced to Text format before assignment
# so Excel keeps the exact textual representation (e.g. trailing zeros) instead
# of silently converting the cell to a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.914.79"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "1.921.43"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.38"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5057"
$ws.Range("E7").Value = "  -2.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4055"
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08357"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.34"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.105"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.87"
$ws.Range("E12").Value = "  +3.60%  "
$ws.Range("D13").Value = "1.923.63"
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.244"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.26"
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001099"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("E19").Value = "  -1.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.29"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "29.959.30"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.35"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("D26").Value = "2.146.71"
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.06"
$ws.Range("E27").Value = "  +2.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.44"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.338"
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.94"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("E31").Value = "  +4.20%  "
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.969"
$ws.Range("E33").Value = "  -2.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.768"
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.424"
$ws.Range("E35").Value = "  +1.57%  "
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06423"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2157"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6598"
$ws.Range("E39").Value = "  +1.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.762"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.196"
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.40"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.213"
$ws.Range("E43").Value = "  -1.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.231"
$ws.Range("E44").Value = "  +8.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.50"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6106"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.617"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "121.87"
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.07"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.127"
$ws.Range("E51").Value = "  -2.84%  "
